$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.643.99"
$ws.Range("D3").Value = "1.887.56"
$ws.Range("E3").Value = "  +1.60%  "
$ws.Range("D5").Value = "'237.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "'0.4835"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("D8").Value = "'0.2859"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.13%  "
$ws.Range("D9").Value = "'0.06542"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.22%  "
$ws.Range("D10").Value = "1.830.96"
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").Value = "'16.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.40%  "
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").Value = "'87.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").Value = "'0.6632"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.79%  "
$ws.Range("D16").Value = "30.602.65"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").Value = "'13.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "'0.000007574"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").Value = "'229.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.48%  "
$ws.Range("D21").Value = "2.103.99"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'5.260"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").Value = "'6.179"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.76%  "
$ws.Range("D25").Value = "'9.410"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.34%  "
$ws.Range("D26").Value = "'167.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.33%  "
$ws.Range("D27").Value = "'18.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("D28").Value = "'1.951"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("D29").Value = "'0.1021"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +11.04%  "
$ws.Range("D30").Value = "'1.395"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.67%  "
$ws.Range("D31").Value = "'4.326"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.07%  "
$ws.Range("D32").Value = "'4.015"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.34%  "
$ws.Range("D33").Value = "'0.05045"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.75%  "
$ws.Range("D34").Value = "'1.199"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.00%  "
$ws.Range("D35").Value = "'0.7475"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.38%  "
$ws.Range("D36").Value = "'0.9999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").Value = "'2.714"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.07%  "
$ws.Range("E38").Value = "  +2.91%  "
$ws.Range("D39").Value = "'2.656"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.25%  "
$ws.Range("D40").Value = "'0.9205"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.51%  "
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("D42").Value = "'107.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("D43").Value = "'0.4267"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").Value = "'5.636"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.54%  "
$ws.Range("D46").Value = "'7.407"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.87%  "
$ws.Range("D47").Value = "'64.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("D48").Value = "'0.1272"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.03%  "
$ws.Range("D49").Value = "'1.476"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("D50").Value = "'8.930"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.69%  "
$ws.Range("D51").Value = "'33.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.61%  "
